# Scheduled refresh of market-board figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# for the Leve-profit tracking sheets. Values sourced from the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9: Distill, My Heart
$ws.Range("H9").Value = 268.5
$ws.Range("I9").Value = 268.5
$ws.Range("K9").Value = 268.5
$ws.Range("M9").Value = -99.5

# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 1082.28
$ws.Range("I15").Value = 1082.28
$ws.Range("K15").Value = 3246.84
$ws.Range("M15").Value = -3077.84

# Row 37: The Wailers' First Law of Potion
$ws.Range("H37").Value = 6673.75
$ws.Range("J37").Value = 7231.6665
$ws.Range("L37").Value = 21694.9995
$ws.Range("N37").Value = -21946.9995

# Row 45: The House Always Wins
$ws.Range("H45").Value = 1508
$ws.Range("I45").Value = 1017
$ws.Range("K45").Value = 3051
$ws.Range("M45").Value = -2859

# Row 48: The Sting of Conscience
$ws.Range("H48").Value = 552.1667
$ws.Range("I48").Value = 483.66666
$ws.Range("J48").Value = 620.6667
$ws.Range("K48").Value = 1450.99998
$ws.Range("L48").Value = 1862.0001
$ws.Range("M48").Value = -1158.99998
$ws.Range("N48").Value = -2446.0001

# Row 49: Going Nowhere Fast
$ws.Range("H49").Value = 25001572
$ws.Range("I49").Value = 140
$ws.Range("K49").Value = 420
$ws.Range("M49").Value = -284

# Row 56: Sleepless in Silvertear
$ws.Range("H56").Value = 552.1667
$ws.Range("I56").Value = 483.66666
$ws.Range("J56").Value = 620.6667
$ws.Range("K56").Value = 1450.99998
$ws.Range("L56").Value = 1862.0001
$ws.Range("M56").Value = -916.9999800000001
$ws.Range("N56").Value = -2930.0001

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 3880.8333
$ws.Range("I76").Value = 3307
$ws.Range("K76").Value = 3307
$ws.Range("M76").Value = -2992

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 3880.8333
$ws.Range("I79").Value = 3307
$ws.Range("K79").Value = 3307
$ws.Range("M79").Value = -2215

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2514.5334
$ws.Range("J138").Value = 2949.5264
$ws.Range("L138").Value = 8848.5792
$ws.Range("N138").Value = -19128.5792


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 195.5
$ws.Range("I5").Value = 113
$ws.Range("K5").Value = 113
$ws.Range("M5").Value = -1

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 35841.56
$ws.Range("I32").Value = 21204.174
$ws.Range("K32").Value = 21204.174
$ws.Range("M32").Value = -20917.174

# Row 37: Get Shirty
$ws.Range("H37").Value = 12519208
$ws.Range("J37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("N37").Value = -30546

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2459.4
$ws.Range("I63").Value = 2459.4
$ws.Range("K63").Value = 2459.4
$ws.Range("M63").Value = -1773.4

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2459.4
$ws.Range("I66").Value = 2459.4
$ws.Range("K66").Value = 12297
$ws.Range("M66").Value = -8865

# Row 80: A Squire to Inspire
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -21996

# Row 83: All's Fair in Highborn Assassination (L)
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -69984

# Row 92: Mail It In
$ws.Range("H92").Value = 30550
$ws.Range("J92").Value = 30550
$ws.Range("L92").Value = 30550
$ws.Range("N92").Value = -35542

# Row 97: Ore for Me
$ws.Range("H97").Value = 12822329
$ws.Range("I97").Value = 15874077
$ws.Range("J97").Value = 4984
$ws.Range("K97").Value = 15874077
$ws.Range("L97").Value = 4984
$ws.Range("M97").Value = -15873581
$ws.Range("N97").Value = -5976

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1977.6666
$ws.Range("I122").Value = 1883.2142
$ws.Range("K122").Value = 5649.642599999999
$ws.Range("M122").Value = -3199.642599999999

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1626.1628
$ws.Range("I132").Value = 1184.6666
$ws.Range("K132").Value = 3553.9998
$ws.Range("M132").Value = -1023.9998


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Range("H4").Value = 195.5
$ws.Range("I4").Value = 113
$ws.Range("K4").Value = 113
$ws.Range("M4").Value = 2

# Row 22: Riveting Run
$ws.Range("H22").Value = 291
$ws.Range("I22").Value = 291
$ws.Range("K22").Value = 291
$ws.Range("M22").Value = -118

# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 15750.667
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 15750.667
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 4849.304
$ws.Range("I105").Value = 4551.75
$ws.Range("K105").Value = 4551.75
$ws.Range("M105").Value = -2804.75

# Row 132: Always Be Prepaired
$ws.Range("H132").Value = 86997
$ws.Range("J132").Value = 86997
$ws.Range("L132").Value = 86997
$ws.Range("N132").Value = -97117


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 436.5
$ws.Range("J5").Value = 425
$ws.Range("L5").Value = 1275
$ws.Range("N5").Value = -1499

# Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 331.25
$ws.Range("J23").Value = 331.25
$ws.Range("L23").Value = 993.75
$ws.Range("N23").Value = -1463.75

# Row 107: Slippery Service
$ws.Range("H107").Value = 866.8570999999999
$ws.Range("I107").Value = 1029
$ws.Range("J107").Value = 776.7778
$ws.Range("K107").Value = 3087
$ws.Range("L107").Value = 2330.3334
$ws.Range("M107").Value = -1167
$ws.Range("N107").Value = -6170.3334

# Row 132: More Mezcal
$ws.Range("H132").Value = 1734.4546
$ws.Range("I132").Value = 1563.1666
$ws.Range("J132").Value = 1940
$ws.Range("K132").Value = 14068.4994
$ws.Range("L132").Value = 17460
$ws.Range("M132").Value = -11538.4994
$ws.Range("N132").Value = -22520

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 436.5
$ws.Range("J135").Value = 425
$ws.Range("L135").Value = 3825
$ws.Range("N135").Value = -8895

# Row 139: Najoothie
$ws.Range("H139").Value = 4053.75
$ws.Range("I139").Value = 4664.3335
$ws.Range("J139").Value = 2222
$ws.Range("K139").Value = 13993.0005
$ws.Range("L139").Value = 6666
$ws.Range("M139").Value = -8853.000499999998
$ws.Range("N139").Value = -16946


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 6060.25
$ws.Range("I80").Value = 3499
$ws.Range("J80").Value = 7597
$ws.Range("K80").Value = 3499
$ws.Range("L80").Value = 7597
$ws.Range("M80").Value = -2501
$ws.Range("N80").Value = -9593

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 6060.25
$ws.Range("I83").Value = 3499
$ws.Range("J83").Value = 7597
$ws.Range("K83").Value = 17495
$ws.Range("L83").Value = 37985
$ws.Range("M83").Value = -12503
$ws.Range("N83").Value = -47969


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1718.25
$ws.Range("I22").Value = 1457.8334
$ws.Range("K22").Value = 1457.8334
$ws.Range("M22").Value = -1162.8334

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1718.25
$ws.Range("I27").Value = 1457.8334
$ws.Range("K27").Value = 1457.8334
$ws.Range("M27").Value = -1350.8334

# Row 40: Best Served Toad
$ws.Range("H40").Value = 3460.2856
$ws.Range("I40").Value = 2084.2856
$ws.Range("J40").Value = 6212.2856
$ws.Range("K40").Value = 2084.2856
$ws.Range("L40").Value = 6212.2856
$ws.Range("M40").Value = -1948.2856
$ws.Range("N40").Value = -6484.2856

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2086.2
$ws.Range("J46").Value = 2014.2858
$ws.Range("L46").Value = 2014.2858
$ws.Range("N46").Value = -2390.2858

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1901.5
$ws.Range("I61").Value = 1901.5
$ws.Range("K61").Value = 1901.5
$ws.Range("M61").Value = -1699.5

# Row 113: Peace in Rest
$ws.Range("H113").Value = 1901.5
$ws.Range("I113").Value = 1901.5
$ws.Range("K113").Value = 1901.5
$ws.Range("M113").Value = 268.5

# Row 133: The Perfect Accessory
$ws.Range("H133").Value = 99992.5
$ws.Range("J133").Value = 99992.5
$ws.Range("L133").Value = 99992.5
$ws.Range("N133").Value = -105052.5


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 24603.69
$ws.Range("I132").Value = 24679.135
$ws.Range("K132").Value = 74037.405
$ws.Range("M132").Value = -71507.405

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 26250
$ws.Range("J136").Value = 2500
$ws.Range("L136").Value = 7500
$ws.Range("N136").Value = -12600
